# Apply the "Förändrad" (changed) date update and add friendly display
# names to the HYPERLINK() formulas, matching the upstream automatic
# export update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update column C ("Förändrad") from 45184 -> 45186 for every data row ---
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row   # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)   # column C
    if ($cCell.Value2 -eq 45184) {
        $cCell.Value = 45186
    }
}

# --- 2. Add the "Beteckning" (column A) as the HYPERLINK() friendly name ---
# Any of columns S, T, U, V, W, X, Y that already hold a HYPERLINK formula
# with a single argument get a second argument appended with the row's
# "Beteckning" text, e.g.
#   HYPERLINK("...A 48258-2020.xlsx")
#   -> HYPERLINK("...A 48258-2020.xlsx", "A 48258-2020")
$hyperlinkCols = 19, 20, 21, 22, 23, 24, 25   # S..Y

for ($r = 2; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Value2   # column A = Beteckning
    if ([string]::IsNullOrEmpty($name)) { continue }

    foreach ($c in $hyperlinkCols) {
        $cell = $ws.Cells.Item($r, $c)
        $f = $cell.Formula
        if ([string]::IsNullOrEmpty($f)) { continue }
        if ($f -notmatch '^=HYPERLINK\(') { continue }
        if ($f -match ',') { continue }   # already has a second argument

        $newFormula = $f.Substring(0, $f.Length - 1) + ', "' + $name + '")'
        # (upstream export keeps a space after the comma, e.g. `..., "Name")`)
        $cell.Formula = $newFormula
    }
}
